# Update automatico via Actualizar 02-05-2021 18-25-56
#
# 1) Refresh the timestamp stored in D716:D729 (last existing block) to a
#    slightly re-computed serial value.
# 2) Append a brand-new 14-row block (rows 730-743) repeating the same
#    Nombre/URL/Disponibilidad/Fecha cycle used throughout the sheet, each
#    row carrying a live hyperlink in column B, and a new shared timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Tiny refresh of the existing last block's Fecha column.
# ---------------------------------------------------------------------
$refreshedDate = 44232.7462469676
for ($r = 716; $r -le 729; $r++) {
    $ws.Cells.Item($r, 4).Value = $refreshedDate
}

# ---------------------------------------------------------------------
# 2) Append the new 14-row block (rows 730-743).
# ---------------------------------------------------------------------
$names = @(
    "Odoo",
    "Blackbox",
    "PowerBI",
    "Dropbox",
    "Odoo",
    "GEE",
    "UtilidadesOdoo",
    "Filtros Dashboard",
    "MapStore",
    "GeoServer",
    "Tomcat",
    "Shiny",
    "Github",
    "EZ Exporter"
)

$urls = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)

# Display text for column B; identical to target URL except MapStore which
# also shows the trailing "#/" fragment (split off into the hyperlink's
# SubAddress/location when the relationship is created).
$displayUrls = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)

$subAddresses = @("", "", "", "", "", "", "", "", "/", "", "", "", "", "")

$newDate = 44232.76784993634
$startRow = 730

for ($i = 0; $i -lt $names.Count; $i++) {
    $r = $startRow + $i

    $ws.Cells.Item($r, 1).Value = $names[$i]

    $bCell = $ws.Cells.Item($r, 2)
    $bCell.Value = $displayUrls[$i]
    if ($subAddresses[$i] -ne "") {
        $ws.Hyperlinks.Add($bCell, $urls[$i], $subAddresses[$i])
    } else {
        $ws.Hyperlinks.Add($bCell, $urls[$i])
    }
    $bCell.Style = "Hyperlink"

    $ws.Cells.Item($r, 3).Value = "Disponible"

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value = $newDate
    $dCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
